$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 34: Sophomore Slump
$ws.Range("H34").Value2 = 7596.4
$ws.Range("I34").Value2 = 7596.4
$ws.Range("K34").Value2 = 7596.4
$ws.Range("M34").Value2 = -7393.4

# Row 36: You Put Your Left Hand In
$ws.Range("H36").Value2 = 7596.4
$ws.Range("I36").Value2 = 7596.4
$ws.Range("K36").Value2 = 7596.4
$ws.Range("M36").Value2 = -6881.4

# Row 62: The Mustache Suits Him
$ws.Range("H62").Value2 = 4719.2
$ws.Range("J62").Value2 = 9148.666999999999
$ws.Range("L62").Value2 = 9148.666999999999
$ws.Range("N62").Value2 = -10396.667

# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value2 = 4719.2
$ws.Range("J65").Value2 = 9148.666999999999
$ws.Range("L65").Value2 = 45743.335
$ws.Range("N65").Value2 = -51983.335

# Row 98: The Dotted Line
$ws.Range("H98").Value2 = 328.33334
$ws.Range("I98").Value2 = 328.33334
$ws.Range("K98").Value2 = 328.33334
$ws.Range("M98").Value2 = 1169.66666

# Row 116: Growing Up
$ws.Range("H116").Value2 = 4005
$ws.Range("I116").Value2 = 4005
$ws.Range("K116").Value2 = 4005
$ws.Range("M116").Value2 = -563

# Row 122: Wishful Inking
$ws.Range("H122").Value2 = 328.33334
$ws.Range("I122").Value2 = 328.33334
$ws.Range("K122").Value2 = 985.0000200000001
$ws.Range("M122").Value2 = 1464.99998

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value2 = 2062.35
$ws.Range("I45").Value2 = 1231.1818
$ws.Range("K45").Value2 = 1231.1818
$ws.Range("M45").Value2 = -854.1818000000001

# Row 59: Parasitic Win
$ws.Range("H59").Value2 = 0
$ws.Range("J59").Value2 = 0
$ws.Range("L59").Value2 = 0
$ws.Range("N59").ClearContents()

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value2 = 6213.769
$ws.Range("I61").Value2 = 3968.4285
$ws.Range("K61").Value2 = 3968.4285
$ws.Range("M61").Value2 = -3756.4285

# Row 97: Ore for Me
$ws.Range("H97").Value2 = 678.0833
$ws.Range("I97").Value2 = 679.7778
$ws.Range("J97").Value2 = 673
$ws.Range("K97").Value2 = 679.7778
$ws.Range("L97").Value2 = 673
$ws.Range("M97").Value2 = -183.7778
$ws.Range("N97").Value2 = -1665

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value2 = 1765.8667
$ws.Range("J132").Value2 = 1499
$ws.Range("L132").Value2 = 4497
$ws.Range("N132").Value2 = -9557

# Row 136: Metal with Mettle
$ws.Range("H136").Value2 = 6213.769
$ws.Range("I136").Value2 = 3968.4285
$ws.Range("K136").Value2 = 11905.2855
$ws.Range("M136").Value2 = -9355.2855

$ws = $wb.Worksheets.Item("BSM")
# Row 40: Can You Spare a Dolabra
$ws.Range("H40").Value2 = 40000
$ws.Range("J40").Value2 = 40000
$ws.Range("L40").Value2 = 40000
$ws.Range("N40").Value2 = -40530

# Row 54: Get Me to the War on Time
$ws.Range("H54").Value2 = 5046.4
$ws.Range("I54").Value2 = 4058
$ws.Range("J54").Value2 = 9000
$ws.Range("K54").Value2 = 4058
$ws.Range("L54").Value2 = 9000
$ws.Range("M54").Value2 = -3574
$ws.Range("N54").Value2 = -9968

# Row 94: High Steal
$ws.Range("H94").Value2 = 507.6
$ws.Range("I94").Value2 = 415.75
$ws.Range("K94").Value2 = 415.75
$ws.Range("M94").Value2 = 35.25

# Row 96: Hammer Time
$ws.Range("H96").Value2 = 14050.4
$ws.Range("I96").Value2 = 14050.4
$ws.Range("J96").Value2 = 0
$ws.Range("K96").Value2 = 14050.4
$ws.Range("L96").Value2 = 0
$ws.Range("M96").Value2 = -11304.4
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value2 = 5880.0605
$ws.Range("I31").Value2 = 3307.7
$ws.Range("J31").Value2 = 9837.538
$ws.Range("K31").Value2 = 3307.7
$ws.Range("L31").Value2 = 9837.538
$ws.Range("M31").Value2 = -3012.7
$ws.Range("N31").Value2 = -10427.538

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value2 = 5880.0605
$ws.Range("I34").Value2 = 3307.7
$ws.Range("J34").Value2 = 9837.538
$ws.Range("K34").Value2 = 3307.7
$ws.Range("L34").Value2 = 9837.538
$ws.Range("M34").Value2 = -3105.7
$ws.Range("N34").Value2 = -10241.538

# Row 94: Beech, Please
$ws.Range("H94").Value2 = 5050.8335
$ws.Range("J94").Value2 = 7300.75
$ws.Range("L94").Value2 = 7300.75
$ws.Range("N94").Value2 = -8202.75

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value2 = 2617.6
$ws.Range("I132").Value2 = 2617.6
$ws.Range("K132").Value2 = 7852.799999999999
$ws.Range("M132").Value2 = -5322.799999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 33: Cooking with Gas
$ws.Range("H33").Value2 = 307.85715
$ws.Range("I33").Value2 = 307.85715
$ws.Range("K33").Value2 = 1847.1429
$ws.Range("M33").Value2 = -1564.1429

# Row 114: One Last Meal
$ws.Range("H114").Value2 = 415.83334
$ws.Range("I114").Value2 = 168.77777
$ws.Range("K114").Value2 = 506.33331
$ws.Range("M114").Value2 = 2747.66669

# Row 120: A Happy End
$ws.Range("H120").Value2 = 4199.3335
$ws.Range("I120").Value2 = 4199.3335
$ws.Range("J120").Value2 = 0
$ws.Range("K120").Value2 = 12598.0005
$ws.Range("L120").Value2 = 0
$ws.Range("M120").Value2 = -7760.000499999998
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 68: Sense of Entitlement
$ws.Range("H68").Value2 = 45258
$ws.Range("J68").Value2 = 45258
$ws.Range("L68").Value2 = 45258
$ws.Range("N68").Value2 = -46880

# Row 71: Charting the Trends (L)
$ws.Range("H71").Value2 = 45258
$ws.Range("J71").Value2 = 45258
$ws.Range("L71").Value2 = 135774
$ws.Range("N71").Value2 = -143886

# Row 86: Keeping Claw and Order
$ws.Range("H86").Value2 = 2000
$ws.Range("J86").Value2 = 2000
$ws.Range("L86").Value2 = 2000
$ws.Range("N86").Value2 = -4372

# Row 89: Ring of Reciprocity (L)
$ws.Range("H89").Value2 = 2000
$ws.Range("J89").Value2 = 2000
$ws.Range("L89").Value2 = 6000
$ws.Range("N89").Value2 = -17856

# Row 97: If I'd a Koppranickel for Every Time...
$ws.Range("H97").Value2 = 360.92856
$ws.Range("I97").Value2 = 330.1
$ws.Range("K97").Value2 = 330.1
$ws.Range("M97").Value2 = 165.9

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value2 = 979.2222
$ws.Range("J22").Value2 = 847
$ws.Range("L22").Value2 = 847
$ws.Range("N22").Value2 = -1437

# Row 27: Fire and Hide
$ws.Range("H27").Value2 = 979.2222
$ws.Range("J27").Value2 = 847
$ws.Range("L27").Value2 = 847
$ws.Range("N27").Value2 = -1061

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value2 = 6999.6
$ws.Range("J68").Value2 = 9999.333000000001
$ws.Range("L68").Value2 = 9999.333000000001
$ws.Range("N68").Value2 = -11497.333

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value2 = 6999.6
$ws.Range("J71").Value2 = 9999.333000000001
$ws.Range("L71").Value2 = 49996.665
$ws.Range("N71").Value2 = -57484.665

# Row 74: Overall, We Blend In
$ws.Range("H74").Value2 = 0
$ws.Range("I74").Value2 = 0
$ws.Range("K74").Value2 = 0
$ws.Range("M74").ClearContents()

# Row 77: Eviction Notice (L)
$ws.Range("H77").Value2 = 0
$ws.Range("I77").Value2 = 0
$ws.Range("K77").Value2 = 0
$ws.Range("M77").ClearContents()

# Row 93: Hide to Go Seek
$ws.Range("H93").Value2 = 646.1667
$ws.Range("I93").Value2 = 646.1667
$ws.Range("K93").Value2 = 646.1667
$ws.Range("M93").Value2 = 601.8333

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope
$ws.Range("H136").Value2 = 3184.52
$ws.Range("I136").Value2 = 2085.1875
$ws.Range("K136").Value2 = 6255.5625
$ws.Range("M136").Value2 = -3705.5625
